$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the neighboring header cell (G1) into the new header cell (H1)
$ws.Range("G1").Copy($ws.Range("H1"))

# Set the header text
$ws.Range("H1").Value = "Save"

# Add the "Save" values for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
